$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The pre-existing A1:B1 carried a bordered/centered style; the appended
# data below is plain, so clear any inherited formatting first.
$ws.Range("A1:E2").ClearFormats()

# Header row
$ws.Range("A1").Value = "MIGRATION DATE"
$ws.Range("B1").Value = "FINANCIAL INSTITUTION NAME"
$ws.Range("C1").Value = "ENTITY ID"
$ws.Range("D1").Value = "ADDRESS"
$ws.Range("E1").Value = "CITY"

# Data row - leading apostrophe keeps the date-looking value as literal
# text instead of letting Excel auto-convert it to a date serial number.
$ws.Range("A2").Value = "'2025-10-16"
$ws.Range("B2").Value = "YYY"
$ws.Range("C2").Value = "123ABX007"
$ws.Range("D2").Value = "Karapakkam"
$ws.Range("E2").Value = "Chennai"
